$wb = $excel.ActiveWorkbook

# --- 1. Rename Sheet1 -> portfolio_input ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Name = "portfolio_input"

# --- 2. Insert a new worksheet right after portfolio_input, named portfolio_metadata ---
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "portfolio_metadata"

# --- 3. Populate the row labels of the portfolio_metadata sheet first ---
$newSheet.Range("A1").Value = "Metadata"
$newSheet.Range("A2").Value = "Name"
$newSheet.Range("A3").Value = "Creator"
$newSheet.Range("A4").Value = "Date"
$newSheet.Range("A5").Value = "Description"

# --- 4. Update the portfolio name / description / creator cells (drop label prefixes) ---
$ws1.Range("B3").Value = "Proportional Portfolio"
$ws1.Range("D3").Value = "The Proportional portfolio reflects a buildout of energy resources in the SJV consistent with the idea that the SJV contributes to the State's energy goals proportional to its resources. That proportionality is defined differently for each feedstock to commodity pathway and should be thought of as a guiding principle rather than a strict rule."
$ws1.Range("G3").Value = "Nidhi Kalra (nidhi@rand.org)"

# --- 5. Link the metadata values back to portfolio_input and set the date ---
$newSheet.Range("B2").Formula = "=portfolio_input!B3"
$newSheet.Range("B3").Formula = "=portfolio_input!D3"

$newSheet.Range("B4").NumberFormat = "mm-dd-yy"
$newSheet.Range("B4").Value = Get-Date -Year 2024 -Month 1 -Day 31 -Hour 0 -Minute 0 -Second 0

$newSheet.Range("B5").Formula = "=portfolio_input!D3"

# --- 6. Add the column header for the value column last ---
$newSheet.Range("B1").Value = "Metadata Value"

# --- 7. Restore selection / view state ---
$newSheet.Range("A1:B5").Select()
$ws1.Range("E10").Select()
$ws1.Activate()
